# Insert a new record row for "Femacal de La Calera - Poroto granado"
# directly below the existing row 90, shifting all subsequent rows
# (old 91..204) down by one to (92..205), then populate the new row
# with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 91; Excel shifts rows 91:204 down to 92:205
# and extends the used range (dimension) automatically.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new observation.
$ws.Range("A91").Value = 3
$ws.Range("B91").Value = "Femacal de La Calera"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44902
$ws.Range("E91").Value = 5
$ws.Range("F91").Value = 100112030
$ws.Range("G91").Value = "Poroto granado"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 68
$ws.Range("K91").Value = 35000
$ws.Range("L91").Value = 36000
$ws.Range("M91").Value = 35559
$ws.Range("N91").Value = "$/saco 25 kilos"
$ws.Range("O91").Value = "Provincia de Limarí"
$ws.Range("P91").Value = 1422
$ws.Range("Q91").Value = 25
$ws.Range("R91").Value = "Hortaliza"
